# Atualizei dados da bibi e add
# Updates column K (new data series) and recalculates column AG (total)
# for rows 2-6 on Sheet1, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => [K new value, AG new value]
$updates = @(
    @{ Row = 2; K = 12113; AG = 87704.73 },
    @{ Row = 3; K = 3766;  AG = 37926.05 },
    @{ Row = 4; K = 3135;  AG = 31754.4 },
    @{ Row = 5; K = 2252;  AG = 29617.05 },
    @{ Row = 6; K = 21266; AG = 187002.23 }
)

foreach ($u in $updates) {
    $ws.Range("K$($u.Row)").Value = $u.K
    $ws.Range("AG$($u.Row)").Value = $u.AG
}

$wb.Save()
